$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct value assignments (safe, non-numeric-looking strings) ---
$ws.Range("D2").Value = "28.814.87"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "1.882.13"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "1.917.39"
$ws.Range("E12").Value = "  +5.22%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "28.825.08"
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "2.120.85"
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E43").Value = "  +4.29%  "
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("E51").Value = "  +0.52%  "

# --- Numeric-looking D values: write as text formulas, then convert formula -> value via PasteSpecial ---
# (keeps cell type as text/string instead of Excel auto-converting to a number)
$ws.Range("D4").Formula = "=""1.006"""
$ws.Range("D5").Formula = "=""323.32"""
$ws.Range("D6").Formula = "=""1.004"""
$ws.Range("D7").Formula = "=""0.4668"""
$ws.Range("D8").Formula = "=""0.3935"""
$ws.Range("D9").Formula = "=""0.07923"""
$ws.Range("D10").Formula = "=""0.9820"""
$ws.Range("D11").Formula = "=""22.35"""
$rng = $ws.Range("D4:D11")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D13").Formula = "=""5.747"""
$ws.Range("D14").Formula = "=""7.014"""
$ws.Range("D15").Formula = "=""0.06975"""
$ws.Range("D16").Formula = "=""88.81"""
$rng = $ws.Range("D13:D16")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D18").Formula = "=""0.00001011"""
$ws.Range("D19").Formula = "=""16.98"""
$ws.Range("D20").Formula = "=""1.003"""
$rng = $ws.Range("D18:D20")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D22").Formula = "=""5.352"""
$ws.Range("D23").Formula = "=""11.10"""
$ws.Range("D24").Formula = "=""2.117"""
$rng = $ws.Range("D22:D24")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D26").Formula = "=""153.69"""
$ws.Range("D27").Formula = "=""19.42"""
$ws.Range("D28").Formula = "=""5.773"""
$ws.Range("D29").Formula = "=""2.004"""
$ws.Range("D30").Formula = "=""120.07"""
$ws.Range("D31").Formula = "=""0.09384"""
$ws.Range("D32").Formula = "=""0.9400"""
$ws.Range("D33").Formula = "=""5.320"""
$ws.Range("D34").Formula = "=""1.362"""
$ws.Range("D35").Formula = "=""3.345"""
$ws.Range("D36").Formula = "=""0.05921"""
$ws.Range("D37").Formula = "=""0.02127"""
$ws.Range("D38").Formula = "=""1.159"""
$rng = $ws.Range("D26:D38")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D40").Formula = "=""0.5734"""
$ws.Range("D41").Formula = "=""0.1799"""
$ws.Range("D42").Formula = "=""9.997"""
$ws.Range("D43").Formula = "=""0.07302"""
$ws.Range("D44").Formula = "=""11.86"""
$ws.Range("D45").Formula = "=""1.176"""
$ws.Range("D46").Formula = "=""0.5353"""
$ws.Range("D47").Formula = "=""2.127"""
$ws.Range("D48").Formula = "=""1.850"""
$ws.Range("D49").Formula = "=""114.33"""
$ws.Range("D50").Formula = "=""2.370"""
$rng = $ws.Range("D40:D50")
$rng.Copy()
$rng.PasteSpecial(-4163)

$excel.CutCopyMode = $false
